# Adds the single-echelon spare-parts inventory parameters:
#   - order_cost / lead_time columns on the existing "spare_parts" sheet
#   - four new sheets: holding_costs, spare_parts_required, max_capacity, reorder_level

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Extend the "spare_parts" sheet with order_cost / lead_time columns
# ---------------------------------------------------------------------------
$spareParts = $wb.Worksheets.Item("spare_parts")

$spareParts.Range("B1").Value = "order_cost"
$spareParts.Range("C1").Value = "lead_time"

$spareParts.Range("B2").Value = 10
$spareParts.Range("C2").Value = 2

$spareParts.Range("B3").Value = 20
$spareParts.Range("C3").Value = 2

$spareParts.Range("B4").Value = 30
$spareParts.Range("C4").Value = 3

$spareParts.Range("B5").Value = 40
$spareParts.Range("C5").Value = 3

# ---------------------------------------------------------------------------
# 2. Add the four new sheets, in order, after "capacity_base_vessels"
# ---------------------------------------------------------------------------
$capacityBaseVessels = $wb.Worksheets.Item("capacity_base_vessels")

$holdingCosts = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $capacityBaseVessels)
$holdingCosts.Name = "holding_costs"

$sparePartsRequired = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $holdingCosts)
$sparePartsRequired.Name = "spare_parts_required"

$maxCapacity = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sparePartsRequired)
$maxCapacity.Name = "max_capacity"

$reorderLevel = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $maxCapacity)
$reorderLevel.Name = "reorder_level"

# ---------------------------------------------------------------------------
# 3. holding_costs: bases (columns) x spare_parts (rows)
# ---------------------------------------------------------------------------
$holdingCosts.Range("B1").Formula = '=TRANSPOSE(FILTER(bases!A2:A100, bases!A2:A100<>""))'
$holdingCosts.Range("A2").Formula = '=FILTER(spare_parts!A2:A100, spare_parts!A2:A100<>"")'

$holdingCosts.Range("B2").Value = 5
$holdingCosts.Range("C2").Value = 8

$holdingCosts.Range("B3").Value = 3
$holdingCosts.Range("C3").Value = 5

$holdingCosts.Range("B4").Value = 4
$holdingCosts.Range("C4").Value = 4

$holdingCosts.Range("B5").Value = 2
$holdingCosts.Range("C5").Value = 1

# ---------------------------------------------------------------------------
# 4. spare_parts_required: tasks (columns) x spare_parts (rows)
# ---------------------------------------------------------------------------
$sparePartsRequired.Range("B1").Formula = '=TRANSPOSE(FILTER(tasks!A2:A100, tasks!A2:A100<>""))'
$sparePartsRequired.Range("A2").Formula = '=FILTER(spare_parts!A2:A100, spare_parts!A2:A100<>"")'

$sparePartsRequired.Range("B2").Value = 2
$sparePartsRequired.Range("C2").Value = 0
$sparePartsRequired.Range("D2").Value = 0

$sparePartsRequired.Range("B3").Value = 1
$sparePartsRequired.Range("C3").Value = 0
$sparePartsRequired.Range("D3").Value = 1

$sparePartsRequired.Range("B4").Value = 0
$sparePartsRequired.Range("C4").Value = 2
$sparePartsRequired.Range("D4").Value = 0

$sparePartsRequired.Range("B5").Value = 1
$sparePartsRequired.Range("C5").Value = 1
$sparePartsRequired.Range("D5").Value = 2

# ---------------------------------------------------------------------------
# 5. max_capacity: bases (columns) x spare_parts (rows)
# ---------------------------------------------------------------------------
$maxCapacity.Range("B1").Formula = '=TRANSPOSE(FILTER(bases!A2:A100, bases!A2:A100<>""))'
$maxCapacity.Range("A2").Formula = '=FILTER(spare_parts!A2:A100, spare_parts!A2:A100<>"")'

$maxCapacity.Range("B2").Value = 15
$maxCapacity.Range("C2").Value = 10

$maxCapacity.Range("B3").Value = 15
$maxCapacity.Range("C3").Value = 10

$maxCapacity.Range("B4").Value = 14
$maxCapacity.Range("C4").Value = 8

$maxCapacity.Range("B5").Value = 14
$maxCapacity.Range("C5").Value = 8

# ---------------------------------------------------------------------------
# 6. reorder_level: bases (columns) x spare_parts (rows)
# ---------------------------------------------------------------------------
$reorderLevel.Range("B1").Formula = '=TRANSPOSE(FILTER(bases!A2:A100, bases!A2:A100<>""))'
$reorderLevel.Range("A2").Formula = '=FILTER(spare_parts!A2:A100, spare_parts!A2:A100<>"")'

$reorderLevel.Range("B2").Value = 5
$reorderLevel.Range("C2").Value = 5

$reorderLevel.Range("B3").Value = 5
$reorderLevel.Range("C3").Value = 5

$reorderLevel.Range("B4").Value = 5
$reorderLevel.Range("C4").Value = 5

$reorderLevel.Range("B5").Value = 5
$reorderLevel.Range("C5").Value = 5

# ---------------------------------------------------------------------------
# 7. Selections / active-cell bookkeeping on touched sheets
# ---------------------------------------------------------------------------
$spareParts.Activate()
$spareParts.Range("G9").Select()

$capacityBaseVessels.Activate()
$capacityBaseVessels.Range("K24").Select()

$holdingCosts.Activate()
$holdingCosts.Range("B1").Select()

$sparePartsRequired.Activate()
$sparePartsRequired.Range("A2").Select()

$reorderLevel.Activate()
$reorderLevel.Range("C7").Select()

# max_capacity ends up as the active tab/sheet, matching the target workbook view
$maxCapacity.Activate()
$maxCapacity.Range("C7").Select()
